$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-10: columns A-D are cluster/symbol strings, E-T are numeric metrics
# This reflects the corrected NATMI LR-pair output (full 3x3 cluster cross-product)

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cthrc1"
$ws.Range("C2").Value = "Fzd6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.9744206666666667
$ws.Range("H2").Value = 2.923262
$ws.Range("I2").Value = 0.0323957296480766
$ws.Range("J2").Value = 0.03239572964807661
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.877148
$ws.Range("N2").Value = 41.631444
$ws.Range("O2").Value = 0.9551716529386821
$ws.Range("P2").Value = 0.9551716529386822
$ws.Range("Q2").Value = 13.522179805592
$ws.Range("R2").Value = 121.699618250328
$ws.Range("S2").Value = 0.030943482636108
$ws.Range("T2").Value = 0.03094348263610801

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cthrc1"
$ws.Range("C3").Value = "Fzd6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.9744206666666667
$ws.Range("H3").Value = 2.923262
$ws.Range("I3").Value = 0.0323957296480766
$ws.Range("J3").Value = 0.03239572964807661
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4259926666666667
$ws.Range("N3").Value = 1.277978
$ws.Range("O3").Value = 0.02932130719941569
$ws.Range("P3").Value = 0.02932130719941569
$ws.Range("Q3").Value = 0.4150960582484445
$ws.Range("R3").Value = 3.735864524236001
$ws.Range("S3").Value = 0.0009498851409604728
$ws.Range("T3").Value = 0.0009498851409604731

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cthrc1"
$ws.Range("C4").Value = "Fzd6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.9744206666666667
$ws.Range("H4").Value = 2.923262
$ws.Range("I4").Value = 0.0323957296480766
$ws.Range("J4").Value = 0.03239572964807661
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.225293
$ws.Range("N4").Value = 0.675879
$ws.Range("O4").Value = 0.01550703986190206
$ws.Range("P4").Value = 0.01550703986190207
$ws.Range("Q4").Value = 0.2195301552553333
$ws.Range("R4").Value = 1.975771397298
$ws.Range("S4").Value = 0.0005023618710081264
$ws.Range("T4").Value = 0.0005023618710081266

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cthrc1"
$ws.Range("C5").Value = "Fzd6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 28.021538
$ws.Range("H5").Value = 84.06461399999999
$ws.Range("I5").Value = 0.9316080830640275
$ws.Range("J5").Value = 0.9316080830640276
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.877148
$ws.Range("N5").Value = 41.631444
$ws.Range("O5").Value = 0.9551716529386821
$ws.Range("P5").Value = 0.9551716529386822
$ws.Range("Q5").Value = 388.859030013624
$ws.Range("R5").Value = 3499.731270122616
$ws.Range("S5").Value = 0.8898456325913042
$ws.Range("T5").Value = 0.8898456325913044

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cthrc1"
$ws.Range("C6").Value = "Fzd6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 28.021538
$ws.Range("H6").Value = 84.06461399999999
$ws.Range("I6").Value = 0.9316080830640275
$ws.Range("J6").Value = 0.9316080830640276
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4259926666666667
$ws.Range("N6").Value = 1.277978
$ws.Range("O6").Value = 0.02932130719941569
$ws.Range("P6").Value = 0.02932130719941569
$ws.Range("Q6").Value = 11.93696969672133
$ws.Range("R6").Value = 107.432727270492
$ws.Range("S6").Value = 0.02731596679297912
$ws.Range("T6").Value = 0.02731596679297913

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cthrc1"
$ws.Range("C7").Value = "Fzd6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 28.021538
$ws.Range("H7").Value = 84.06461399999999
$ws.Range("I7").Value = 0.9316080830640275
$ws.Range("J7").Value = 0.9316080830640276
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.225293
$ws.Range("N7").Value = 0.675879
$ws.Range("O7").Value = 0.01550703986190206
$ws.Range("P7").Value = 0.01550703986190207
$ws.Range("Q7").Value = 6.313056360633999
$ws.Range("R7").Value = 56.817507245706
$ws.Range("S7").Value = 0.01444648367974404
$ws.Range("T7").Value = 0.01444648367974405

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cthrc1"
$ws.Range("C8").Value = "Fzd6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.082717666666667
$ws.Range("H8").Value = 3.248153
$ws.Range("I8").Value = 0.03599618728789584
$ws.Range("J8").Value = 0.03599618728789584
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 13.877148
$ws.Range("N8").Value = 41.631444
$ws.Range("O8").Value = 0.9551716529386821
$ws.Range("P8").Value = 0.9551716529386822
$ws.Range("Q8").Value = 15.025033302548
$ws.Range("R8").Value = 135.225299722932
$ws.Range("S8").Value = 0.03438253771126985
$ws.Range("T8").Value = 0.03438253771126985

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cthrc1"
$ws.Range("C9").Value = "Fzd6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.082717666666667
$ws.Range("H9").Value = 3.248153
$ws.Range("I9").Value = 0.03599618728789584
$ws.Range("J9").Value = 0.03599618728789584
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4259926666666667
$ws.Range("N9").Value = 1.277978
$ws.Range("O9").Value = 0.02932130719941569
$ws.Range("P9").Value = 0.02932130719941569
$ws.Range("Q9").Value = 0.4612297860704445
$ws.Range("R9").Value = 4.151068074634
$ws.Range("S9").Value = 0.001055455265476096
$ws.Range("T9").Value = 0.001055455265476096

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cthrc1"
$ws.Range("C10").Value = "Fzd6"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.082717666666667
$ws.Range("H10").Value = 3.248153
$ws.Range("I10").Value = 0.03599618728789584
$ws.Range("J10").Value = 0.03599618728789584
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.225293
$ws.Range("N10").Value = 0.675879
$ws.Range("O10").Value = 0.01550703986190206
$ws.Range("P10").Value = 0.01550703986190207
$ws.Range("Q10").Value = 0.2439287112763333
$ws.Range("R10").Value = 2.195358401487
$ws.Range("S10").Value = 0.0005581943111498931
$ws.Range("T10").Value = 0.0005581943111498932
